# Material classification.docx edit script
# 1. Remove the entire Slovenian section (and the blank/"***" separator
#    paragraphs) that used to precede the English text, leaving the
#    English "In the Excel workbook..." paragraph as the new first
#    paragraph of the document.
# 2. Update several English paragraphs' wording.
# 3. Merge the "info.txt" paragraph into an expanded "info.json"
#    paragraph that now also folds in the content that used to live in
#    the standalone "RT properties.txt" / "Ranges.txt" / "Fields.txt"
#    paragraphs, then delete those three now-redundant paragraphs.

$d = $word.ActiveDocument

# --- 1. Delete the leading Slovenian block -------------------------------
# That block is the set of paragraphs before the paragraph that starts
# with "In the Excel workbook". Locate it and delete everything before it.
$cutParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "In the Excel workbook*") {
        $cutParaIndex = $i
        break
    }
}
if ($cutParaIndex -gt 1) {
    $delStart = $d.Paragraphs(1).Range.Start
    $delEnd = $d.Paragraphs($cutParaIndex - 1).Range.End
    $d.Range($delStart, $delEnd).Delete()
}

$find = $d.Content.Find

# --- 2a. "Flags" paragraph: wording unchanged, just normalize runs -------
$find.ClearFormatting()
$find.Execute(
    "- Flags. These fields serve to indicate in the code to which group the material belongs in terms of functionality within the heating system. If the material is invariant, it cannot be anything else. Otherwise, it can have several flags.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Flags. These fields serve to indicate in the code to which group the material belongs in terms of functionality within the heating system. If the material is invariant, it cannot be anything else. Otherwise, it can have several flags.",
    2)

# --- 2b. "Data on thermal properties" paragraph: wording changes ---------
$find.ClearFormatting()
$find.Execute(
    "- Data on thermal properties marked in red must always be present in the material folder. These are the properties at room temperature (~293 K) that must be available in case we do not know the temperature dependencies. These boolean values are therefore always 1.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Data on thermal properties marked in red are the properties at room temperature (~293 K) that should be available in case we do not know the temperature dependencies. These boolean values are therefore always 1 (can also be 0 if we really do not have the data).",
    2)

# --- 2c. "Other data" paragraph: wording unchanged, just normalize runs --
$find.ClearFormatting()
$find.Execute(
    "- Other data are linked to the group within the heating system. Invariant materials, for example, have files rho(T), cp(T), k(T), but magnetocaloric materials don" + [char]0x2019 + "t have these, because there will be separate cp(T), namely for heating and cooling at several magnetic fields (see blue colored fields). There could be even more combinations of listed properties and dependencies, and they depend on the material.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Other data are linked to the group within the heating system. Invariant materials, for example, have files rho(T), cp(T), k(T), but magnetocaloric materials don" + [char]0x2019 + "t have these, because there will be separate cp(T), namely for heating and cooling at several magnetic fields (see blue colored fields). There could be even more combinations of listed properties and dependencies, and they depend on the material.",
    2)

# --- 2d. "info.txt" -> expanded "info.json" paragraph --------------------
$find.ClearFormatting()
$find.Execute(
    "- The info.txt file contains the material ID, short name and melting temperature in sequence, separated by commas.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- The info.json file is in JSON format and must contain: RT properties, ranges, fields, ID, short name, long name, and melting point. The RT properties must contain the values of density, specific heat, conductivity and emissivity at room temperature. The ranges must contain the temperature ranges in which each material property is defined. Temperature ranges must be recorded for density, specific heat capacity, thermal conductivity, adiabatic temperature change (if the material is caloric), emissivity, and other relevant properties (e.g. seebeck coefficient for thermoelectric materials, etc.). The fields contains strengths of the external fields where properties for caloric materials are defined. When it comes to magnetic fields, the values are in T, when it comes to electric fields, they are in MV/m, and for pressure and stress, they are in bars. (See any material for example).",
    2)

# --- 3. Delete the now-redundant "RT properties.txt" / "Ranges.txt" /
#        "Fields.txt" paragraphs (their content now lives in info.json) --
$targets = @(
    "- The RT properties.txt file is in JSON format and must contain the values of density, specific heat, conductivity and emissivity at room temperature.",
    "- The Ranges.txt file is in JSON format and must contain the temperature ranges in which each material property is defined. Temperature ranges must be recorded for density, specific heat capacity, thermal conductivity, adiabatic temperature change (if the material is caloric), emissivity, and other relevant properties (e.g. seebeck coefficient for thermoelectric materials, etc.).",
    "- The Fields.txt file contains a column with the strengths of the external fields where properties for caloric materials are defined. When it comes to magnetic fields, the values are in T, when it comes to electric fields, they are in MV/m, and for pressure and stress, they are in bars."
)

foreach ($target in $targets) {
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text.TrimEnd("`r", "`a") -eq $target) {
            $p.Range.Delete()
            break
        }
    }
}
